# Apply the "exp.id" / "pot" ID-variable addition to the cfu_bacteria sheet.
# Adds two new columns:
#   D "pot"    - replicate/pot number, cycling 1,2,3 down the data rows
#   E "exp.id" - concatenation formula: =A&"_"&B&"_"&"bact"&"_"&D

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Headers (row 1) ------------------------------------------------------
$ws.Range("D1").Value = "pot"
$ws.Range("E1").Value = "exp.id"

# --- D column: pot number, 1,2,3 repeating for rows 2..88 -----------------
$potValues = New-Object 'object[,]' 87,1
for ($i = 0; $i -lt 87; $i++) {
    $potValues[$i,0] = ($i % 3) + 1
}
$ws.Range("D2:D88").Value = $potValues

# --- E column: exp.id formula ----------------------------------------------
# Entered the same way a user would: first cell on its own, then two
# separate fills (the data has a repeated header row at row 67 which is
# where the author's second fill operation started), reproducing the
# shared-formula groupings seen in the saved file.
$ws.Range("E2").Formula = "=A2&""_""&B2&""_""&""bact""&""_""&D2"
$ws.Range("E3:E66").Formula = "=A3&""_""&B3&""_""&""bact""&""_""&D3"
$ws.Range("E67:E88").Formula = "=A67&""_""&B67&""_""&""bact""&""_""&D67"

# --- Selection / view state -------------------------------------------------
$ws.Range("A2:E88").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 61
